# Auto-generated script applying scheduled market-data refresh updates
# to the Leve profit calculation sheets (currentAveragePrice / Leve price /
# Leve profit columns H-N) across all character sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 8120
$ws.Range("I64").Value = 11100
$ws.Range("K64").Value = 11100
$ws.Range("M64").Value = -10852
$ws.Range("H67").Value = 8120
$ws.Range("I67").Value = 11100
$ws.Range("K67").Value = 11100
$ws.Range("M67").Value = -10242
$ws.Range("H76").Value = 62502760
$ws.Range("I76").Value = 71431304
$ws.Range("J76").Value = 2966.6667
$ws.Range("K76").Value = 71431304
$ws.Range("L76").Value = 2966.6667
$ws.Range("M76").Value = -71430989
$ws.Range("N76").Value = -3596.6667
$ws.Range("H79").Value = 62502760
$ws.Range("I79").Value = 71431304
$ws.Range("J79").Value = 2966.6667
$ws.Range("K79").Value = 71431304
$ws.Range("L79").Value = 2966.6667
$ws.Range("M79").Value = -71430212
$ws.Range("N79").Value = -5150.6667
$ws.Range("H116").Value = 4634.4
$ws.Range("I116").Value = 4666.643
$ws.Range("K116").Value = 4666.643
$ws.Range("M116").Value = -1224.643
$ws.Range("H124").Value = 36200
$ws.Range("J124").Value = 36200
$ws.Range("L124").Value = 36200
$ws.Range("N124").Value = -46020
$ws.Range("H133").Value = 66221.53999999999
$ws.Range("J133").Value = 72807.27
$ws.Range("L133").Value = 72807.27
$ws.Range("N133").Value = -82927.27
$ws.Range("H134").Value = 50780
$ws.Range("J134").Value = 50780
$ws.Range("L134").Value = 50780
$ws.Range("N134").Value = -60920
$ws.Range("H136").Value = 83400
$ws.Range("I136").Value = 77000
$ws.Range("J136").Value = 89800
$ws.Range("K136").Value = 77000
$ws.Range("L136").Value = 89800
$ws.Range("M136").Value = -71900
$ws.Range("N136").Value = -100000
$ws.Range("H138").Value = 3986.3635
$ws.Range("I138").Value = 2655.1765
$ws.Range("J138").Value = 5400.75
$ws.Range("K138").Value = 7965.529500000001
$ws.Range("L138").Value = 16202.25
$ws.Range("M138").Value = -2825.529500000001
$ws.Range("N138").Value = -26482.25
$ws.Range("H139").Value = 78000
$ws.Range("J139").Value = 78000
$ws.Range("L139").Value = 78000
$ws.Range("N139").Value = -88280
$ws.Range("H140").Value = 58276
$ws.Range("J140").Value = 58276
$ws.Range("L140").Value = 58276
$ws.Range("N140").Value = -68636

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 45915.78
$ws.Range("I2").Value = 57443.055
$ws.Range("J2").Value = 4417.6
$ws.Range("K2").Value = 57443.055
$ws.Range("L2").Value = 4417.6
$ws.Range("M2").Value = -57330.055
$ws.Range("N2").Value = -4643.6
$ws.Range("H32").Value = 28369.021
$ws.Range("I32").Value = 18827.057
$ws.Range("J32").Value = 58729.816
$ws.Range("K32").Value = 18827.057
$ws.Range("L32").Value = 58729.816
$ws.Range("M32").Value = -18540.057
$ws.Range("N32").Value = -59303.816
$ws.Range("H61").Value = 1306.45
$ws.Range("I61").Value = 1296.2632
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1296.2632
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -1084.2632
$ws.Range("N61").Value = -1924
$ws.Range("H63").Value = 3034.1458
$ws.Range("I63").Value = 2210.2727
$ws.Range("K63").Value = 2210.2727
$ws.Range("M63").Value = -1524.2727
$ws.Range("H66").Value = 3034.1458
$ws.Range("I66").Value = 2210.2727
$ws.Range("K66").Value = 11051.3635
$ws.Range("M66").Value = -7619.363499999999
$ws.Range("H74").Value = 1817.963
$ws.Range("I74").Value = 1862.421
$ws.Range("J74").Value = 1712.375
$ws.Range("K74").Value = 1862.421
$ws.Range("L74").Value = 1712.375
$ws.Range("M74").Value = -988.421
$ws.Range("N74").Value = -3460.375
$ws.Range("H77").Value = 1817.963
$ws.Range("I77").Value = 1862.421
$ws.Range("J77").Value = 1712.375
$ws.Range("K77").Value = 9312.105
$ws.Range("L77").Value = 8561.875
$ws.Range("M77").Value = -4944.105
$ws.Range("N77").Value = -17297.875
$ws.Range("H88").Value = 2454.0715
$ws.Range("I88").Value = 2608.3333
$ws.Range("J88").Value = 2338.375
$ws.Range("K88").Value = 2608.3333
$ws.Range("L88").Value = 2338.375
$ws.Range("M88").Value = -2202.3333
$ws.Range("N88").Value = -3150.375
$ws.Range("H91").Value = 2454.0715
$ws.Range("I91").Value = 2608.3333
$ws.Range("J91").Value = 2338.375
$ws.Range("K91").Value = 2608.3333
$ws.Range("L91").Value = 2338.375
$ws.Range("M91").Value = -1204.3333
$ws.Range("N91").Value = -5146.375
$ws.Range("H116").Value = 45915.78
$ws.Range("I116").Value = 57443.055
$ws.Range("J116").Value = 4417.6
$ws.Range("K116").Value = 57443.055
$ws.Range("L116").Value = 4417.6
$ws.Range("M116").Value = -55149.055
$ws.Range("N116").Value = -9005.6
$ws.Range("H122").Value = 1036
$ws.Range("I122").Value = 1036
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3108
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -658
$ws.Range("H132").Value = 4159.7144
$ws.Range("I132").Value = 3531
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 10593
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -8063
$ws.Range("N132").Value = -20054
$ws.Range("H136").Value = 1306.45
$ws.Range("I136").Value = 1296.2632
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 3888.7896
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -1338.7896
$ws.Range("N136").Value = -9600
$ws.Range("M122").Value = ""

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 45915.78
$ws.Range("I3").Value = 57443.055
$ws.Range("J3").Value = 4417.6
$ws.Range("K3").Value = 57443.055
$ws.Range("L3").Value = 4417.6
$ws.Range("M3").Value = -57329.055
$ws.Range("N3").Value = -4645.6
$ws.Range("H134").Value = 39705.98
$ws.Range("I134").Value = 1881.7174
$ws.Range("J134").Value = 288265.44
$ws.Range("K134").Value = 5645.1522
$ws.Range("L134").Value = 864796.3200000001
$ws.Range("M134").Value = -3110.1522
$ws.Range("N134").Value = -869866.3200000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1890.25
$ws.Range("I16").Value = 1207.3334
$ws.Range("J16").Value = 2300
$ws.Range("K16").Value = 1207.3334
$ws.Range("L16").Value = 2300
$ws.Range("M16").Value = -920.3334
$ws.Range("N16").Value = -2874
$ws.Range("H62").Value = 6945.5557
$ws.Range("I62").Value = 6786.2856
$ws.Range("K62").Value = 6786.2856
$ws.Range("M62").Value = -6162.2856
$ws.Range("H65").Value = 6945.5557
$ws.Range("I65").Value = 6786.2856
$ws.Range("K65").Value = 33931.428
$ws.Range("M65").Value = -30811.428
$ws.Range("H113").Value = 1890.25
$ws.Range("I113").Value = 1207.3334
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 1207.3334
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 962.6666
$ws.Range("N113").Value = -6640
$ws.Range("H122").Value = 1887.2142
$ws.Range("I122").Value = 1961.3334
$ws.Range("J122").Value = 1442.5
$ws.Range("K122").Value = 5884.0002
$ws.Range("L122").Value = 4327.5
$ws.Range("M122").Value = -3434.0002
$ws.Range("N122").Value = -9227.5
$ws.Range("H132").Value = 1544.6471
$ws.Range("I132").Value = 1145.5834
$ws.Range("J132").Value = 2502.4
$ws.Range("K132").Value = 3436.7502
$ws.Range("L132").Value = 7507.200000000001
$ws.Range("M132").Value = -906.7501999999999
$ws.Range("N132").Value = -12567.2

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4045
$ws.Range("I137").Value = 4045
$ws.Range("K137").Value = 12135
$ws.Range("M137").Value = -7035

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17843.242
$ws.Range("I70").Value = 34461.152
$ws.Range("K70").Value = 34461.152
$ws.Range("M70").Value = -34191.152
$ws.Range("H73").Value = 17843.242
$ws.Range("I73").Value = 34461.152
$ws.Range("K73").Value = 34461.152
$ws.Range("M73").Value = -33525.152
$ws.Range("H80").Value = 3575.4827
$ws.Range("I80").Value = 4023.4614
$ws.Range("K80").Value = 4023.4614
$ws.Range("M80").Value = -3025.4614
$ws.Range("H83").Value = 3575.4827
$ws.Range("I83").Value = 4023.4614
$ws.Range("K83").Value = 20117.307
$ws.Range("M83").Value = -15125.307
$ws.Range("H107").Value = 341.54544
$ws.Range("I107").Value = 203.38889
$ws.Range("J107").Value = 963.25
$ws.Range("K107").Value = 203.38889
$ws.Range("L107").Value = 963.25
$ws.Range("M107").Value = 1716.61111
$ws.Range("N107").Value = -4803.25
$ws.Range("H113").Value = 2842.3635
$ws.Range("I113").Value = 1435.6
$ws.Range("J113").Value = 4014.6667
$ws.Range("K113").Value = 1435.6
$ws.Range("L113").Value = 4014.6667
$ws.Range("M113").Value = 734.4000000000001
$ws.Range("N113").Value = -8354.6667
$ws.Range("H126").Value = 4838.727
$ws.Range("I126").Value = 3479.111
$ws.Range("J126").Value = 10957
$ws.Range("K126").Value = 10437.333
$ws.Range("L126").Value = 32871
$ws.Range("M126").Value = -7967.332999999999
$ws.Range("N126").Value = -37811
$ws.Range("H132").Value = 3410.3784
$ws.Range("I132").Value = 2974.4062
$ws.Range("J132").Value = 6200.6
$ws.Range("K132").Value = 8923.2186
$ws.Range("L132").Value = 18601.8
$ws.Range("M132").Value = -6393.2186
$ws.Range("N132").Value = -23661.8

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1050
$ws.Range("I46").Value = 1050
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1050
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = -862
$ws.Range("H122").Value = 2895.647
$ws.Range("I122").Value = 2674.8333
$ws.Range("J122").Value = 3144.0625
$ws.Range("K122").Value = 8024.499899999999
$ws.Range("L122").Value = 9432.1875
$ws.Range("M122").Value = -5574.499899999999
$ws.Range("N122").Value = -14332.1875
$ws.Range("H136").Value = 5491.1797
$ws.Range("I136").Value = 3568.2727
$ws.Range("J136").Value = 6246.607
$ws.Range("K136").Value = 10704.8181
$ws.Range("L136").Value = 18739.821
$ws.Range("M136").Value = -8154.8181
$ws.Range("N136").Value = -23839.821
$ws.Range("M46").Value = ""

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("H113").Value = 608.7
$ws.Range("I113").Value = 390.33334
$ws.Range("J113").Value = 936.25
$ws.Range("K113").Value = 1171.00002
$ws.Range("L113").Value = 2808.75
$ws.Range("M113").Value = 998.9999800000001
$ws.Range("N113").Value = -7148.75
$ws.Range("N10").Value = ""
